# DEI-3-3 Added data model sketch
# Fills in the "Real hours" column (B) for the existing backlog rows and
# adds a new "Tech debt" mini-section (rows 16-20) below the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B3: "TODO" flagged in bold red (set font color/bold BEFORE the value so the
# engine reuses/creates the bold-red font the same way Excel would).
$ws.Range("B3").Font.Color = 255
$ws.Range("B3").Font.Bold = $true
$ws.Range("B3").Value = "TODO"

# --- New "Tech debt" block under the table (rows 16-20) ---------------------

$ws.Range("A16").Value = "DEI-0"
$ws.Range("B16").Value = "Tech debt:"
$ws.Range("B17").Value = "udf"
$ws.Range("B18").Value = "apply appropriate names (classes, utils, tests)"

# --- Column B "Real hours" estimates for the existing rows ------------------

$ws.Range("B4").Value = "50 min"
$ws.Range("B5").Value = "5 min"
$ws.Range("B8").Value = "5 min"

$ws.Range("B19").Value = "parametrise queries impliying big data"

$ws.Range("B6").Value = "10 min"
$ws.Range("B7").Value = "10 min"

# B9 is a plain number (20), styled with the existing red-font style.
$ws.Range("B9").Font.Color = 255
$ws.Range("B9").Value = 20

# B20 stays empty but picks up a red, time-formatted (h:mm) style.
$ws.Range("B20").Font.Color = 255
$ws.Range("B20").NumberFormat = "h:mm"

# --- View state: scroll down a row and leave the selection on B10 ----------

$excel.ActiveWindow.ScrollRow = 2
$ws.Range("B10").Select()
